# Estadisticos Matutinos 15 Oct
# Fills in the "Reprobados" / "Aprobados" counts and the resulting
# percentages + averages for the 1st partial / 2nd partial / Final
# statistics sheets, and lists the students who are "rescatable"
# (candidates to retake) on the Rescatables sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Estadisticos 1P
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")

# Row 2 - 5AEV
$ws1.Cells.Item(2,4).Value = 5        # D Blancos
$ws1.Cells.Item(2,6).Value = 17       # F Aprobados
$ws1.Cells.Item(2,7).Value = 77.27    # G Por_Apro
$ws1.Cells.Item(2,8).Value = 6.9      # H Promedio

# Row 3 - 5ALCV
$ws1.Cells.Item(3,4).Value = 4
$ws1.Cells.Item(3,6).Value = 27
$ws1.Cells.Item(3,7).Value = 87.1
$ws1.Cells.Item(3,8).Value = 8.5

# Row 4 - 5APV
$ws1.Cells.Item(4,4).Value = 14
$ws1.Cells.Item(4,5).Value = 2        # E Reprobados
$ws1.Cells.Item(4,6).Value = 22
$ws1.Cells.Item(4,7).Value = 57.89
$ws1.Cells.Item(4,8).Value = 8.4

# Row 5 - 5ARHV
$ws1.Cells.Item(5,4).Value = 10
$ws1.Cells.Item(5,5).Value = 1
$ws1.Cells.Item(5,6).Value = 24
$ws1.Cells.Item(5,7).Value = 68.57
$ws1.Cells.Item(5,8).Value = 7.9

# Row 6 - 5ASV
$ws1.Cells.Item(6,4).Value = 9
$ws1.Cells.Item(6,5).Value = 3
$ws1.Cells.Item(6,6).Value = 21
$ws1.Cells.Item(6,7).Value = 63.64
$ws1.Cells.Item(6,8).Value = 7.8

# ---------------------------------------------------------------
# Estadisticos 2P (only Reprobados changes here)
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")

$ws2.Cells.Item(2,5).Value = 17   # 5AEV
$ws2.Cells.Item(3,5).Value = 27   # 5ALCV
$ws2.Cells.Item(4,5).Value = 24   # 5APV
$ws2.Cells.Item(5,5).Value = 25   # 5ARHV
$ws2.Cells.Item(6,5).Value = 24   # 5ASV

# ---------------------------------------------------------------
# Estadisticos Final (same values as Estadisticos 1P)
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Estadisticos Final")

$ws3.Cells.Item(2,4).Value = 5
$ws3.Cells.Item(2,6).Value = 17
$ws3.Cells.Item(2,7).Value = 77.27
$ws3.Cells.Item(2,8).Value = 6.9

$ws3.Cells.Item(3,4).Value = 4
$ws3.Cells.Item(3,6).Value = 27
$ws3.Cells.Item(3,7).Value = 87.1
$ws3.Cells.Item(3,8).Value = 8.5

$ws3.Cells.Item(4,4).Value = 14
$ws3.Cells.Item(4,5).Value = 2
$ws3.Cells.Item(4,6).Value = 22
$ws3.Cells.Item(4,7).Value = 57.89
$ws3.Cells.Item(4,8).Value = 8.4

$ws3.Cells.Item(5,4).Value = 10
$ws3.Cells.Item(5,5).Value = 1
$ws3.Cells.Item(5,6).Value = 24
$ws3.Cells.Item(5,7).Value = 68.57
$ws3.Cells.Item(5,8).Value = 7.9

$ws3.Cells.Item(6,4).Value = 9
$ws3.Cells.Item(6,5).Value = 3
$ws3.Cells.Item(6,6).Value = 21
$ws3.Cells.Item(6,7).Value = 63.64
$ws3.Cells.Item(6,8).Value = 7.8

# ---------------------------------------------------------------
# Rescatables - add the 4 rescatable students
# Columns: NC | Paterno | Materno | Nombres | Nombre_Largo | Grupo | Reprobadas
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Rescatables")

$ws4.Cells.Item(2,1).Value = 19330051920284
$ws4.Cells.Item(2,2).Value = "LOPEZ"
$ws4.Cells.Item(2,3).Value = "ROQUE"
$ws4.Cells.Item(2,4).Value = "CONCEPCION"
$ws4.Cells.Item(2,5).Value = "CÁLCULO INTEGRAL"
$ws4.Cells.Item(2,6).Value = "5ALCV"
$ws4.Cells.Item(2,7).Value = 6

$ws4.Cells.Item(3,1).Value = 19330051920319
$ws4.Cells.Item(3,2).Value = "CRUZ"
$ws4.Cells.Item(3,3).Value = "AULIS"
$ws4.Cells.Item(3,4).Value = "ALDO"
$ws4.Cells.Item(3,5).Value = "CÁLCULO INTEGRAL"
$ws4.Cells.Item(3,6).Value = "5APV"
$ws4.Cells.Item(3,7).Value = 6

$ws4.Cells.Item(4,1).Value = 19330051920433
$ws4.Cells.Item(4,2).Value = "GONZALEZ"
$ws4.Cells.Item(4,3).Value = "SERRANO"
$ws4.Cells.Item(4,4).Value = "CRISTIAN"
$ws4.Cells.Item(4,5).Value = "CÁLCULO INTEGRAL"
$ws4.Cells.Item(4,6).Value = "5ASV"
$ws4.Cells.Item(4,7).Value = 6

$ws4.Cells.Item(5,1).Value = 19330051920253
$ws4.Cells.Item(5,2).Value = "GOMEZ"
$ws4.Cells.Item(5,3).Value = "RIVERA"
$ws4.Cells.Item(5,4).Value = "LUIS ANGEL"
$ws4.Cells.Item(5,5).Value = "CÁLCULO INTEGRAL"
$ws4.Cells.Item(5,6).Value = "5ASV"
$ws4.Cells.Item(5,7).Value = 6
